$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.934.15"
$ws.Range("E2").Value = "'  -0.62%  "

$ws.Range("D3").Value = "'2.550.67"
$ws.Range("E3").Value = "'  -0.28%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.07%  "

$ws.Range("D5").Value = "'304.11"
$ws.Range("E5").Value = "'  +1.24%  "

$ws.Range("D6").Value = "'98.28"
$ws.Range("E6").Value = "'  +5.59%  "

$ws.Range("E7").Value = "'  -0.18%  "

$ws.Range("E8").Value = "'  +0.11%  "

$ws.Range("D9").Value = "'0.545"
$ws.Range("E9").Value = "'  -1.04%  "

$ws.Range("D10").Value = "'36.60"
$ws.Range("E10").Value = "'  +1.68%  "

$ws.Range("D11").Value = "'0.0828"
$ws.Range("E11").Value = "'  +2.24%  "

$ws.Range("B12").Value = "'Polkadot"
$ws.Range("C12").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'7.73"
$ws.Range("E12").Value = "'  +0.17%  "

$ws.Range("B13").Value = "'TRON"
$ws.Range("C13").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.115"
$ws.Range("E13").Value = "'  +5.67%  "

$ws.Range("D14").Value = "'2.942.87"
$ws.Range("E14").Value = "'  -0.14%  "

$ws.Range("D15").Value = "'2.534.20"
$ws.Range("E15").Value = "'  -0.82%  "

$ws.Range("D16").Value = "'14.91"
$ws.Range("E16").Value = "'  +5.23%  "

$ws.Range("E17").Value = "'  +0.27%  "

$ws.Range("D18").Value = "'43.004.87"
$ws.Range("E18").Value = "'  -0.48%  "

$ws.Range("D19").Value = "'13.50"
$ws.Range("E19").Value = "'  +1.18%  "

$ws.Range("D20").Value = "'0.0₃0995"
$ws.Range("E20").Value = "'  +1.29%  "

$ws.Range("D21").Value = "'6.61"
$ws.Range("E21").Value = "'  -0.86%  "

$ws.Range("D22").Value = "'71.94"
$ws.Range("E22").Value = "'  -0.39%  "

$ws.Range("D23").Value = "'253.91"
$ws.Range("E23").Value = "'  -2.70%  "

$ws.Range("D25").Value = "'2.08"
$ws.Range("E25").Value = "'  -2.64%  "

$ws.Range("D26").Value = "'27.94"
$ws.Range("E26").Value = "'  -5.95%  "

$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "'  -0.17%  "

$ws.Range("D28").Value = "'10.14"
$ws.Range("E28").Value = "'  +0.91%  "

$ws.Range("D29").Value = "'37.78"
$ws.Range("E29").Value = "'  +0.49%  "

$ws.Range("E30").Value = "'  -0.93%  "

$ws.Range("D31").Value = "'6.05"
$ws.Range("E31").Value = "'  +0.24%  "

$ws.Range("D32").Value = "'158.45"
$ws.Range("E32").Value = "'  +2.62%  "

$ws.Range("D33").Value = "'2.75"
$ws.Range("E33").Value = "'  -0.64%  "

$ws.Range("E34").Value = "'  -1.34%  "

$ws.Range("D35").Value = "'0.0805"
$ws.Range("E35").Value = "'  +0.47%  "

$ws.Range("B36").Value = "'Celestia"
$ws.Range("C36").Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").Value = "'19.04"
$ws.Range("E36").Value = "'  +12.21%  "

$ws.Range("B37").Value = "'LidoDAOToken"
$ws.Range("C37").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'3.30"
$ws.Range("E37").Value = "'  -2.92%  "

$ws.Range("D38").Value = "'25.95"
$ws.Range("E38").Value = "'  +10.56%  "

$ws.Range("D39").Value = "'0.116"
$ws.Range("E39").Value = "'  -0.70%  "

$ws.Range("E40").Value = "'  -0.41%  "

$ws.Range("D41").Value = "'2.10"
$ws.Range("E41").Value = "'  +32.86%  "

$ws.Range("D42").Value = "'3.44"
$ws.Range("E42").Value = "'  -1.14%  "

$ws.Range("E43").Value = "'  -0.40%  "

$ws.Range("D44").Value = "'2.085.85"
$ws.Range("E44").Value = "'  +0.11%  "

$ws.Range("E45").Value = "'  -2.63%  "

$ws.Range("D47").Value = "'86.53"
$ws.Range("E47").Value = "'  +0.59%  "

$ws.Range("D48").Value = "'8.98"
$ws.Range("E48").Value = "'  +1.50%  "

$ws.Range("D49").Value = "'2.800.32"
$ws.Range("E49").Value = "'  -0.08%  "

$ws.Range("D50").Value = "'74.93"
$ws.Range("E50").Value = "'  +7.75%  "

$ws.Range("D51").Value = "'103.50"
$ws.Range("E51").Value = "'  -1.28%  "

